# Add team record (Wins/Losses/Ties) columns to the KCR_2002 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): AD1=Wins, AE1=Losses, AF1=Ties ---
# Copy the formatting of the existing header cell (AC1 - bold, centered,
# thin border) onto the three new header cells so they match the rest
# of row 1, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (rows 2-54): same W/L/T record for every player ---
$lastRow = 54
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 62   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 100  # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
